$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values (meanEMG / legmaxROM data) for columns B:E
$ws.Range("B2").Value = 35.262305819999995
$ws.Range("C2").Value = 20.995868953125068
$ws.Range("D2").Value = 30.18661815750005
$ws.Range("E2").Value = 28.206329338124988

# Update row 3 values for columns B:E
$ws.Range("B3").Value = 30.923174999999901
$ws.Range("C3").Value = 18.391727160000016
$ws.Range("D3").Value = 30.161621840000009
$ws.Range("E3").Value = 29.790700447499944

# Update the selected range to reflect the updated region
$ws.Range("B1:E3").Select()
